$d = $word.ActiveDocument

$replacements = @(
    @{old="773×4="; new="445×8="},
    @{old="698×4="; new="238×8="},
    @{old="747×5="; new="880×3="},
    @{old="649×5="; new="901×5="},
    @{old="558×5="; new="297×3="},
    @{old="889×3="; new="507×7="},
    @{old="585×4="; new="842×6="},
    @{old="514×2="; new="694×7="},
    @{old="641×7="; new="289×7="},
    @{old="419×7="; new="331×9="},
    @{old="978×4="; new="804×3="},
    @{old="855×4="; new="287×3="},
    @{old="990×6="; new="569×9="},
    @{old="766×6="; new="955×5="},
    @{old="545×4="; new="558×5="},
    @{old="354×5="; new="324×4="},
    @{old="329×9="; new="276×8="},
    @{old="130×2="; new="448×5="},
    @{old="110×3="; new="875×7="},
    @{old="760×5="; new="395×9="},
    @{old="107×6="; new="411×6="},
    @{old="751×4="; new="847×6="},
    @{old="136×7="; new="797×4="},
    @{old="514×3="; new="148×2="},
    @{old="543×8="; new="962×2="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 1)
}
